# Update cryptocurrency price/volume data per the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.437.82"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "1.628.48"
$ws.Range("E3").Value = "  +2.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9955"
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.26"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9971"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3787"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "53.36"
$ws.Range("E8").Value = "  +5.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3660"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.279"
$ws.Range("E10").Value = "  +5.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08204"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9955"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.24"
$ws.Range("E13").Value = "  +6.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.664"
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.463"
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001263"
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("D17").Value = "1.625.33"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.72"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06932"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.40"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.583"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9980"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.99"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").Value = "23.459.54"
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.146"
$ws.Range("E25").Value = "  +13.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.418"
$ws.Range("E26").Value = "  +2.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.38"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.68"
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.283"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.31"
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.413"
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.922"
$ws.Range("E32").Value = "  +5.90%  "
$ws.Range("D33").Value = "1.801.46"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9779"
$ws.Range("E34").Value = "  +3.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02804"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.49"
$ws.Range("E36").Value = "  +4.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07451"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.244"
$ws.Range("E38").Value = "  +3.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2534"
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08846"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.410"
$ws.Range("E41").Value = "  +5.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7175"
$ws.Range("E42").Value = "  +4.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.80"
$ws.Range("E43").Value = "  +6.94%  "
$ws.Range("E44").Value = "  +9.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6622"
$ws.Range("E45").Value = "  +3.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.364"
$ws.Range("E46").Value = "  +5.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.035"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08016"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.43"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.217"
$ws.Range("E51").Value = "  +2.13%  "

Write-Output "Updated cryptos list"
